# Sample Project / Main.xlsx - "SAVE" edit
# Row 11 (the "R40" rule row) had its Rule-name cell (B11) retyped from
# the text "R40" to the text "1". Prefixing the literal with a leading
# apostrophe forces Excel to store the numeric-looking input "1" as text
# (shared string) rather than auto-converting it to a number, matching
# the original text cell B11 was before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").Value = "'1"
